## RC-Flat-Roofs.docx: replace the "Solar panels for high rise buildings"
## inline picture with a hyperlink run pointing at the image's URL on
## ura.gov.sg (the picture's relationship becomes unused; a fresh external
## hyperlink relationship is created and referenced instead).

$d = $word.ActiveDocument

# Locate the inline picture (there is exactly one InlineShape in this
# document - the solar-panel illustration under "Maximum 5m above roof
# level").
$shape = $d.InlineShapes.Item(1)
$rng = $shape.Range
$insertAt = $rng.Start

# Remove the picture run entirely.
$shape.Delete()

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/GFA/GFA54_Solar_Panels_Buildings.jpg?h=100%25&w=100%25"

# Re-insert, in its place, a run-level hyperlink whose visible text is the
# image URL, styled with the built-in "Hyperlink" character style - this
# mirrors the other external hyperlinks already present in the document
# (e.g. the "URA SPACE" / "here" links using rId20/rId27/rId28).
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body><w:p><w:hyperlink r:id="rIdSolarPanelsImgLink"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t xml:space="preserve">__URL_TEXT__</w:t></w:r></w:hyperlink></w:p></w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rIdSolarPanelsImgLink" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="__URL_ATTR__" TargetMode="External"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

# The pkg:xmlData island is unescaped once for element text (so a literal
# "&" must be written as "&amp;" to survive as "&"), but relationship
# attribute values are copied through as-is (so a literal "&" must be
# written bare to come out as a single "&" - writing "&amp;" there would
# double-escape into "&amp;amp;"). Build both spellings from the same
# logical URL.
$urlForText = $url.Replace("&", "&amp;")
$urlForAttr = $url
$xml = $xml.Replace("__URL_TEXT__", $urlForText)
$xml = $xml.Replace("__URL_ATTR__", $urlForAttr)

$ins = $d.Range($insertAt, $insertAt)
$ins.InsertXML($xml)
